$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C1: was an empty numeric cell, now holds the literal text date value "2025-03-05"
# (it must stay text, not get auto-converted by Excel into a date serial number).
# Temporarily format the cell as Text so the string is stored verbatim...
$ws.Range("C1").NumberFormat = "@"
$ws.Range("C1").Value = "2025-03-05"
# ...then copy the original cell formatting back from B1 (same style as before the
# edit: s="1", General number format) so only the value/type of C1 changes.
$ws.Range("B1").Copy() | Out-Null
$ws.Range("C1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# C2: numeric value changes from 3 to 6
$ws.Range("C2").Value = 6
